# Test-Data.xlsx edits:
#  - Search sheet: make date-picking data-driven -> bump CheckIn/CheckOut
#    dates forward by the same 92-day span (new stay: 2026-01-01 .. 2026-01-14)
#    and move the active selection from B2 to C2.
#  - Reservation sheet: move the active selection from B4 (out of range) to A2.

$wb = $excel.ActiveWorkbook

# --- Search sheet -----------------------------------------------------
$search = $wb.Worksheets.Item("Search")
$search.Activate()

$search.Range("B2").Value = Get-Date -Year 2026 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$search.Range("C2").Value = Get-Date -Year 2026 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0

$search.Range("C2").Select()

# --- Reservation sheet --------------------------------------------------
$reservation = $wb.Worksheets.Item("Reservation")
$reservation.Activate()

$reservation.Range("A2").Select()
